# Odds-base update for "Denmark Superligaen" — 08-05-2024 20:15
# Swaps the two fixtures that were recorded with transposed ids/teams/odds
# (rows 178 & 179), corrects the odds for rows 177/178/179/180/181/182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- id column (B) for rows 178/179 swapped -------------------------------
# Force text storage (column is General-formatted, so a bare numeric-looking
# string would otherwise be coerced to a number) before writing the swapped
# ids back, keeping them as text just like the original cells.
$idRange = $ws.Range("B178:B179")
$idRange.NumberFormat = "@"
$ws.Range("B178").Value = "7984029"
$ws.Range("B179").Value = "7984028"

# --- HomeTeam/AwayTeam (E/F) for rows 178/179 swapped ----------------------
$ws.Range("E178").Value = "Viborg"
$ws.Range("F178").Value = "Vejle"
$ws.Range("E179").Value = "Randers FC"
$ws.Range("F179").Value = "Hvidovre IF"

# --- Odds corrections -------------------------------------------------------

# Row 177
$ws.Range("M177").Value = 2.05
$ws.Range("O177").Value = 3.5
$ws.Range("P177").Value = -0.25
$ws.Range("Q177").Value = 1.82
$ws.Range("R177").Value = 2.08
$ws.Range("T177").Value = 1.85
$ws.Range("U177").Value = 2

# Row 178
$ws.Range("J178").Value = 1.75
$ws.Range("K178").Value = 3.6
$ws.Range("L178").Value = 4.333
$ws.Range("M178").Value = 1.8
$ws.Range("N178").Value = 3.5
$ws.Range("O178").Value = 4.2
$ws.Range("P178").Value = -0.5
$ws.Range("Q178").Value = 1.86
$ws.Range("R178").Value = 2.04
$ws.Range("S178").Value = 2.5
$ws.Range("T178").Value = 1.85
$ws.Range("U178").Value = 2

# Row 179
$ws.Range("J179").Value = 1.45
$ws.Range("K179").Value = 4.5
$ws.Range("L179").Value = 5.25
$ws.Range("M179").Value = 1.5
$ws.Range("N179").Value = 4.75
$ws.Range("O179").Value = 5.5
$ws.Range("P179").Value = -1.25
$ws.Range("Q179").Value = 2.05
$ws.Range("R179").Value = 1.85
$ws.Range("S179").Value = 3
$ws.Range("T179").Value = 2
$ws.Range("U179").Value = 1.85

# Row 180
$ws.Range("M180").Value = 3.2
$ws.Range("Q180").Value = 1.95
$ws.Range("R180").Value = 1.95

# Row 181
$ws.Range("M181").Value = 1.5
$ws.Range("O181").Value = 6
$ws.Range("P181").Value = -1
$ws.Range("Q181").Value = 1.81
$ws.Range("R181").Value = 2.09

# Row 182
$ws.Range("M182").Value = 1.363
$ws.Range("N182").Value = 4.75
$ws.Range("O182").Value = 7.5
$ws.Range("Q182").Value = 2.05
$ws.Range("R182").Value = 1.85
$ws.Range("T182").Value = 1.925
$ws.Range("U182").Value = 1.925
